$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("A1").Interior.ThemeColor = 2
$ws.Range("A1").Interior.TintAndShade = -0.0999786
